# Updates cryptos list price (D) / 1h-volume-change (E) columns to
# match the latest scrape, per the commit's row-by-row diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''59.769.74'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.36%  '
$ws.Range("D3").Value = '''2.420.68'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.99%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''551.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.63%  '
$ws.Range("D6").Value = '''137.18'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.73%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  +2.91%  '
$ws.Range("D9").Value = '''0.105'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.55%  '
$ws.Range("E10").Value = '  +2.09%  '
$ws.Range("E11").Value = '  -2.11%  '
$ws.Range("E12").Value = '  +0.55%  '
$ws.Range("D13").Value = '''24.83'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.28%  '
$ws.Range("D14").Value = '''2.851.17'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.01%  '
$ws.Range("D15").Value = '''59.743.56'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.41%  '
$ws.Range("D17").Value = '''2.432.94'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.23%  '
$ws.Range("D18").Value = '''11.27'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.28%  '
$ws.Range("D19").Value = '''4.39'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.14%  '
$ws.Range("D20").Value = '''330.48'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.42%  '
$ws.Range("D21").Value = '''6.67'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.76%  '
$ws.Range("E22").Value = '  -0.07%  '
$ws.Range("D23").Value = '''65.71'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.38%  '
$ws.Range("D24").Value = '''0.172'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.77%  '
$ws.Range("D25").Value = '''8.69'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +5.88%  '
$ws.Range("E26").Value = '  +0.11%  '
$ws.Range("D27").Value = '''1.37'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.69%  '
$sub3 = [char]8323
$ws.Range("D28").Value = "'0.0{0}0777" -f $sub3
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.76%  '
$ws.Range("E29").Value = '  +0.66%  '
$ws.Range("D30").Value = '''170.40'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.20%  '
$ws.Range("D31").Value = '''6.16'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.77%  '
$ws.Range("E32").Value = '  +1.67%  '
$ws.Range("E33").Value = '  +0.70%  '
$ws.Range("D34").Value = '''0.999'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.00%  '
$ws.Range("E35").Value = '  +4.79%  '
$ws.Range("E36").Value = '  +0.06%  '
$ws.Range("D37").Value = '''4.21'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.74%  '
$ws.Range("D38").Value = '''1.61'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.25%  '
$ws.Range("D39").Value = '''39.59'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.95%  '
$ws.Range("D40").Value = '''0.412'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.30%  '
$ws.Range("D41").Value = '''314.78'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +8.73%  '
$ws.Range("D42").Value = '''3.67'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.70%  '
$ws.Range("D43").Value = '''138.52'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.14%  '
$ws.Range("E44").Value = '  +1.68%  '
$ws.Range("E45").Value = '  +1.07%  '
$ws.Range("D46").Value = '''19.53'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.62%  '
$ws.Range("D47").Value = '''0.580'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.55%  '
$ws.Range("D48").Value = '''0.406'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.74%  '
$ws.Range("E49").Value = '  +1.67%  '
$ws.Range("E50").Value = '  +1.30%  '
$ws.Range("E51").Value = '  -0.22%  '
